# #5: property aircraft done
# Fix property_category values that were incorrectly left as "land" for
# rows that actually describe buildings (建物 sheet) and a car (汽車 sheet).

$wb = $excel.ActiveWorkbook

# 建物 (Building) sheet: column I = property_category, rows 2-7 currently "land" -> "building"
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I7").Value = "building"

# 汽車 (Car) sheet: column H = property_category, row 2 currently "land" -> "car"
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
